# Adds a new "disease induced mortality" parameter (columns L/M) to both the
# PerDay and PerSeason sheets, and makes "PerDay" the selected/active sheet
# (previously "PerSeason" was active).
#
# PerDay (sheet1) stores the per-day instantaneous rate, derived from the
# per-season probabilities via -LOG(1-p)/days.
# PerSeason (sheet2) stores the raw per-season probabilities directly.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # PerDay
$ws2 = $wb.Worksheets.Item(2)   # PerSeason

# ---------------------------------------------------------------------------
# PerDay sheet: fill column L first (top to bottom), then column M, so that
# new shared-string entries are created in the same order the workbook
# expects ("disease induced mortality", "???", "McMahon 2023",
# "Carey et al 2006", "Padgett Flohr 2008").
# ---------------------------------------------------------------------------
$ws1.Range("L1").Value = "disease induced mortality"
$ws1.Range("L2").Formula = "=-LOG(1-0.13)/14"
$ws1.Range("L3").Value = "???"
$ws1.Range("L4").Value = "???"
$ws1.Range("L5").Formula = "=-LOG(0.5)/42"
$ws1.Range("L6").Value = 0
$ws1.Range("L7").Value = 0

$ws1.Range("M1").Value = "Citation"
$ws1.Range("M2").Value = "McMahon 2023"
$ws1.Range("M5").Value = "Carey et al 2006"
$ws1.Range("M6").Value = "Daszak 2004"
$ws1.Range("M7").Value = "Padgett Flohr 2008"

# ---------------------------------------------------------------------------
# PerSeason sheet: same layout, but the raw (per-season) probabilities are
# entered directly rather than as a derived formula.
# ---------------------------------------------------------------------------
$ws2.Range("L1").Value = "disease induced mortality"
$ws2.Range("L2").Value = 0.13
$ws2.Range("L3").Value = "???"
$ws2.Range("L4").Value = "???"
$ws2.Range("L5").Value = 1
$ws2.Range("L6").Value = 0
$ws2.Range("L7").Value = 0

$ws2.Range("M1").Value = "Citation"
$ws2.Range("M2").Value = "McMahon 2023"
$ws2.Range("M5").Value = "Carey et al 2006"
$ws2.Range("M6").Value = "Daszak 2004"
$ws2.Range("M7").Value = "Padgett Flohr 2008"

# ---------------------------------------------------------------------------
# New-column widths, matching the workbook's existing "best fit" columns
# (closest values reachable through the ColumnWidth property).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(11).ColumnWidth = 12.65
$ws1.Columns.Item(12).ColumnWidth = 23.3
$ws2.Columns.Item(11).ColumnWidth = 12.65
$ws2.Columns.Item(12).ColumnWidth = 23.3

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the saved view state:
#  - PerSeason's selection moves to the new L:M columns.
#  - PerDay becomes the selected / active sheet (previously PerSeason was).
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("L1:M1048576").Select() | Out-Null

$ws1.Activate()
